# "Generate Report for Handback" - localization-status.xlsx update
#
# This script reflects the handback of the de-de (and completion of the
# zh-cn) localization round: the shared "Status" text flips from
# "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on the per-language sheets get populated, a
# hyperlink is added on the newly-filled "Latest Target File" cell, and
# a few columns are widened so the new long file names are readable.

$wb = $excel.ActiveWorkbook

$targetMd  = "abf1ff20-febe-4748-bfb2-66b3e1478b0d.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d838e510c6d22e494e398ed97bf6089e9ebd3f4c/e2e/abf1ff20-febe-4748-bfb2-66b3e1478b0d.md"
$zhXlf     = "abf1ff20-febe-4748-bfb2-66b3e1478b0d.54ad1d06a6b0dd538de3b4949ff8986929e56b4b.zh-cn.xlf"
$deXlf     = "abf1ff20-febe-4748-bfb2-66b3e1478b0d.54ad1d06a6b0dd538de3b4949ff8986929e56b4b.de-de.xlf"

$handedBackStatus = "Handed back: in sync with en-US"
$zhHandbackTime   = "2016-08-24 21:01:30"
$deHandbackTime   = "2016-08-24 21:01:38"

# Column widths (Excel COM "characters" units). 39.17 / 29.17 are the
# values this host resolves to the workbook's existing 40-char and
# ~30-char stored column widths, matching the other wide columns
# already in the sheets (e.g. "Source File Name"/"Source Path").
$wideColWidth   = 39.17
$mediumColWidth = 29.17

# ---- Overview sheet ------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $handedBackStatus           # zh-cn status
$ov.Range("F2").Value = $handedBackStatus           # de-de status
$ov.Columns.Item(5).ColumnWidth = $mediumColWidth   # E: zh-cn
$ov.Columns.Item(6).ColumnWidth = $mediumColWidth   # F: de-de

# ---- zh-cn sheet -----------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $handedBackStatus           # Status
$zh.Range("I2").Value = $targetMd                   # Latest Target File
$zh.Range("J2").Value = $zhXlf                       # Latest Handback File
$zh.Range("K2").Value = $zhHandbackTime              # Latest Handback DateTime

$zh.Hyperlinks.Add($zh.Range("I2"), $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetMd) | Out-Null

$zh.Columns.Item(3).ColumnWidth = $mediumColWidth    # C: Status
$zh.Columns.Item(9).ColumnWidth = $wideColWidth      # I: Latest Target File
$zh.Columns.Item(10).ColumnWidth = $wideColWidth     # J: Latest Handback File

# ---- de-de sheet -----------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $handedBackStatus           # Status
$de.Range("I2").Value = $targetMd                   # Latest Target File
$de.Range("J2").Value = $deXlf                       # Latest Handback File
$de.Range("K2").Value = $deHandbackTime              # Latest Handback DateTime

$de.Hyperlinks.Add($de.Range("I2"), $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetMd) | Out-Null

$de.Columns.Item(3).ColumnWidth = $mediumColWidth    # C: Status
$de.Columns.Item(9).ColumnWidth = $wideColWidth      # I: Latest Target File
$de.Columns.Item(10).ColumnWidth = $wideColWidth     # J: Latest Handback File

Write-Host "Handback report generated."
